# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with new rates ---
$hoja1 = $wb.Worksheets.Item("Hoja1")

$hoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.75 = 51447.15 pesos`n✅ 51447.15 pesos = 12.66 = 961.88 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update rate table values ---
$tasas = $wb.Worksheets.Item("tasas")

$tasas.Range("N10").Value = 78.43000000000001
$tasas.Range("O10").Value = 4035
$tasas.Range("N12").Value = 4065
$tasas.Range("O12").Value = 76.001
